# Finish setting up effective tax rates: remove the duplicated CRI row
# (row 77, which repeated the COL 2022 trust_wages row already present at
# row 76) from the "trust" sheet. Deleting the whole row shifts every row
# below it up by one, which is what the target workbook shows (CRI/2001
# moves from row 78 to row 77, etc., and the last row 137 disappears).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("trust")

$ws.Rows.Item(77).Delete()
